$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): reorder/insert columns -----------------------
# New column D "PLAN ESTUDIO" is introduced, and the 3 existing headers
# are shifted to A:C ("MATRICULA","NOMBRE","CURP"). Writing D1 first
# keeps the shared-string table ordered the same way the workbook was
# actually authored (PLAN ESTUDIO, MATRICULA, NOMBRE, CURP, ...).
$ws.Range("D1").Value = "PLAN ESTUDIO"
$ws.Range("A1").Value = "MATRICULA"
$ws.Range("B1").Value = "NOMBRE"
$ws.Range("C1").Value = "CURP"

# --- Example data row (row 2) ------------------------------------------
$ws.Range("A2").Value = "206275"
$ws.Range("B2").Value = "LUIS ALEJANDRO"
$ws.Range("C2").Value = "FEVC000117HSRLLRA5"
$ws.Range("D2").Value = "1"

# --- Header formatting: bold the new header row -------------------------
$ws.Range("A1:D1").Font.Bold = $true

# --- Column widths for the new/adjusted columns --------------------------
$ws.Columns.Item(2).ColumnWidth = 17.7109375
$ws.Columns.Item(3).ColumnWidth = 21.85546875
$ws.Columns.Item(4).ColumnWidth = 15.42578125

# --- Selection: header row is selected, matching the saved view ---------
$ws.Range("A1:D1").Select()
